$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H62").Value = 5056.1113
$ws.Range("I62").Value = 4875.625
$ws.Range("J62").Value = 6500
$ws.Range("K62").Value = 4875.625
$ws.Range("L62").Value = 6500
$ws.Range("M62").Value = -4251.625
$ws.Range("N62").Value = -7748
$ws.Range("H65").Value = 5056.1113
$ws.Range("I65").Value = 4875.625
$ws.Range("J65").Value = 6500
$ws.Range("K65").Value = 24378.125
$ws.Range("L65").Value = 32500
$ws.Range("M65").Value = -21258.125
$ws.Range("N65").Value = -38740
$ws.Range("H76").Value = 6073.5
$ws.Range("I76").Value = 5149.75
$ws.Range("J76").Value = 6997.25
$ws.Range("K76").Value = 5149.75
$ws.Range("L76").Value = 6997.25
$ws.Range("M76").Value = -4834.75
$ws.Range("N76").Value = -7627.25
$ws.Range("H79").Value = 6073.5
$ws.Range("I79").Value = 5149.75
$ws.Range("J79").Value = 6997.25
$ws.Range("K79").Value = 5149.75
$ws.Range("L79").Value = 6997.25
$ws.Range("M79").Value = -4057.75
$ws.Range("N79").Value = -9181.25
$ws.Range("H82").Value = 1942.5714
$ws.Range("I82").Value = 1599.6
$ws.Range("J82").Value = 2800
$ws.Range("K82").Value = 4798.799999999999
$ws.Range("L82").Value = 8400
$ws.Range("M82").Value = -4392.799999999999
$ws.Range("N82").Value = -9212
$ws.Range("H85").Value = 1942.5714
$ws.Range("I85").Value = 1599.6
$ws.Range("J85").Value = 2800
$ws.Range("K85").Value = 4798.799999999999
$ws.Range("L85").Value = 8400
$ws.Range("M85").Value = -3394.799999999999
$ws.Range("N85").Value = -11208
$ws.Range("H93").Value = 46866
$ws.Range("J93").Value = 46866
$ws.Range("L93").Value = 46866
$ws.Range("N93").Value = -51858
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2277.4443
$ws.Range("I63").Value = 1449.5
$ws.Range("J63").Value = 3933.3333
$ws.Range("K63").Value = 1449.5
$ws.Range("L63").Value = 3933.3333
$ws.Range("M63").Value = -763.5
$ws.Range("N63").Value = -5305.3333
$ws.Range("H66").Value = 2277.4443
$ws.Range("I66").Value = 1449.5
$ws.Range("J66").Value = 3933.3333
$ws.Range("K66").Value = 7247.5
$ws.Range("L66").Value = 19666.6665
$ws.Range("M66").Value = -3815.5
$ws.Range("N66").Value = -26530.6665
$ws.Range("H88").Value = 7688.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 7688.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 7688.5
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -8500.5
$ws.Range("H91").Value = 7688.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 7688.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 7688.5
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -10496.5
$ws.Range("H102").Value = 1480
$ws.Range("I102").Value = 1480
$ws.Range("K102").Value = 1480
$ws.Range("M102").Value = 142
$ws.Range("H110").Value = 2386.2
$ws.Range("I110").Value = 643.75
$ws.Range("J110").Value = 4999.875
$ws.Range("K110").Value = 643.75
$ws.Range("L110").Value = 4999.875
$ws.Range("M110").Value = 1401.25
$ws.Range("N110").Value = -9089.875
$ws.Range("H132").Value = 2013.64
$ws.Range("I132").Value = 1833.6818
$ws.Range("J132").Value = 3333.3333
$ws.Range("K132").Value = 5501.0454
$ws.Range("L132").Value = 9999.999899999999
$ws.Range("M132").Value = -2971.0454
$ws.Range("N132").Value = -15059.9999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3086.6667
$ws.Range("J20").Value = 3086.6667
$ws.Range("L20").Value = 3086.6667
$ws.Range("N20").Value = -3580.6667
$ws.Range("H86").Value = 3659.8
$ws.Range("J86").Value = 4266.3335
$ws.Range("L86").Value = 4266.3335
$ws.Range("N86").Value = -6512.3335
$ws.Range("H89").Value = 3659.8
$ws.Range("J89").Value = 4266.3335
$ws.Range("L89").Value = 21331.6675
$ws.Range("N89").Value = -32563.6675
$ws.Range("H94").Value = 1188.5416
$ws.Range("I94").Value = 1131.5217
$ws.Range("J94").Value = 2500
$ws.Range("K94").Value = 1131.5217
$ws.Range("L94").Value = 2500
$ws.Range("M94").Value = -680.5217
$ws.Range("N94").Value = -3402
$ws.Range("H99").Value = 2269.7144
$ws.Range("I99").Value = 2314.6667
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 2314.6667
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -816.6667000000002
$ws.Range("N99").Value = -4996
$ws.Range("H100").Value = 33997.57
$ws.Range("J100").Value = 33997.57
$ws.Range("L100").Value = 33997.57
$ws.Range("N100").Value = -36161.57
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 551.8
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 639.75
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 639.75
$ws.Range("M22").Value = 150
$ws.Range("N22").Value = -1339.75
$ws.Range("H31").Value = 4198.2
$ws.Range("I31").Value = 3881.7856
$ws.Range("K31").Value = 3881.7856
$ws.Range("M31").Value = -3586.7856
$ws.Range("H34").Value = 4198.2
$ws.Range("I34").Value = 3881.7856
$ws.Range("K34").Value = 3881.7856
$ws.Range("M34").Value = -3679.7856
$ws.Range("H58").Value = 5495.5
$ws.Range("I58").Value = 3995
$ws.Range("J58").Value = 5995.6665
$ws.Range("K58").Value = 3995
$ws.Range("L58").Value = 5995.6665
$ws.Range("M58").Value = -3792
$ws.Range("N58").Value = -6401.6665
$ws.Range("H86").Value = 4786.375
$ws.Range("I86").Value = 4755.857
$ws.Range("K86").Value = 4755.857
$ws.Range("M86").Value = -3632.857
$ws.Range("H89").Value = 4786.375
$ws.Range("I89").Value = 4755.857
$ws.Range("K89").Value = 23779.285
$ws.Range("M89").Value = -18163.285
$ws.Range("H132").Value = 1819.9
$ws.Range("I132").Value = 1819.9
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5459.700000000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2929.700000000001
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 4622
$ws.Range("I134").Value = 2996.6667
$ws.Range("J134").Value = 9498
$ws.Range("K134").Value = 8990.000100000001
$ws.Range("L134").Value = 28494
$ws.Range("M134").Value = -6455.000100000001
$ws.Range("N134").Value = -33564
$ws.Range("H136").Value = 5495.5
$ws.Range("I136").Value = 3995
$ws.Range("J136").Value = 5995.6665
$ws.Range("K136").Value = 11985
$ws.Range("L136").Value = 17986.9995
$ws.Range("M136").Value = -9435
$ws.Range("N136").Value = -23086.9995
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 264513.47
$ws.Range("I4").Value = 418668.47
$ws.Range("K4").Value = 1256005.41
$ws.Range("M4").Value = -1255893.41
$ws.Range("H131").Value = 2238.25
$ws.Range("I131").Value = 1346.2858
$ws.Range("J131").Value = 2932
$ws.Range("K131").Value = 4038.8574
$ws.Range("L131").Value = 8796
$ws.Range("M131").Value = 1001.1426
$ws.Range("N131").Value = -18876
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 10000
$ws.Range("J39").Value = 10000
$ws.Range("L39").Value = 10000
$ws.Range("N39").Value = -11064
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 4271.25
$ws.Range("I80").Value = 1835
$ws.Range("J80").Value = 5083.3335
$ws.Range("K80").Value = 1835
$ws.Range("L80").Value = 5083.3335
$ws.Range("M80").Value = -837
$ws.Range("N80").Value = -7079.3335
$ws.Range("H83").Value = 4271.25
$ws.Range("I83").Value = 1835
$ws.Range("J83").Value = 5083.3335
$ws.Range("K83").Value = 9175
$ws.Range("L83").Value = 25416.6675
$ws.Range("M83").Value = -4183
$ws.Range("N83").Value = -35400.6675
$ws.Range("H122").Value = 125000000
$ws.Range("I122").Value = 125000000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 375000000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -374997550
$ws.Range("N122").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 899.4
$ws.Range("I22").Value = 899.4
$ws.Range("K22").Value = 899.4
$ws.Range("M22").Value = -604.4
$ws.Range("H27").Value = 899.4
$ws.Range("I27").Value = 899.4
$ws.Range("K27").Value = 899.4
$ws.Range("M27").Value = -792.4
$ws.Range("H46").Value = 2249
$ws.Range("I46").Value = 500
$ws.Range("K46").Value = 500
$ws.Range("M46").Value = -312
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 15000
$ws.Range("J34").Value = 15000
$ws.Range("L34").Value = 15000
$ws.Range("N34").Value = -15406
$ws.Range("H62").Value = 3000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 3000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H68").Value = 54999
$ws.Range("J68").Value = 54999
$ws.Range("L68").Value = 54999
$ws.Range("N68").Value = -56621
$ws.Range("H71").Value = 54999
$ws.Range("J71").Value = 54999
$ws.Range("L71").Value = 164997
$ws.Range("N71").Value = -173109
$ws.Range("H81").Value = 353.93332
$ws.Range("I81").Value = 336.35715
$ws.Range("J81").Value = 600
$ws.Range("K81").Value = 672.7143
$ws.Range("L81").Value = 1200
$ws.Range("M81").Value = 388.2857
$ws.Range("N81").Value = -3322
$ws.Range("H84").Value = 353.93332
$ws.Range("I84").Value = 336.35715
$ws.Range("J84").Value = 600
$ws.Range("K84").Value = 3363.5715
$ws.Range("L84").Value = 6000
$ws.Range("M84").Value = 1940.4285
$ws.Range("N84").Value = -16608
